$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("usually", 1520000000),
    @("complete", 3680000000),
    @("shoot", 567000000),
    @("halls", 116000000),
    @("technology", 5980000000),
    @("hand", 4040000000),
    @("aspect", 495000000),
    @("leave", 2450000000),
    @("science", 2680000000),
    @("echo", 395000000),
    @("tesla", 295000000),
    @("aspect", 495000000),
    @("longer", 1380000000),
    @("simple", 3880000000),
    @("halls", 116000000),
    @("boat", 1560000000),
    @("question", 2780000000),
    @("parent", 629000000),
    @("older", 1610000000),
    @("open", 8690000000),
    @("play", 10110000000),
    @("tesla", 356000000),
    @("erupt", 15100000),
    @("science", 3120000000),
    @("tell", 1940000000),
    @("usually", 1520000000),
    @("soon", 2400000000),
    @("textbook", 197000000),
    @("echo", 323000000),
    @("longer", 1380000000),
    @("open", 8990000000),
    @("shoot", 683000000),
    @("tesla", 1190000000),
    @("science", 2890000000),
    @("time", 15060000000),
    @("make", 9540000000),
    @("older", 1550000000),
    @("limit", 1280000000),
    @("tell", 2020000000),
    @("textbook", 224000000),
    @("longer", 1430000000),
    @("question", 4470000000),
    @("room", 4800000000),
    @("tesla", 313000000),
    @("science", 2820000000),
    @("cheer", 185000000),
    @("limit", 1160000000),
    @("complete", 3140000000),
    @("years", 7360000000),
    @("hand", 4980000000),
    @("time", 13200000000),
    @("afraid", 303000000),
    @("perfect", 3870000000),
    @("tell", 2020000000),
    @("tesla", 313000000),
    @("parent", 660000000),
    @("cultural", 954000000),
    @("time", 13070000000),
    @("older", 1320000000),
    @("open", 7620000000),
    @("halls", 129000000),
    @("simple", 3720000000),
    @("room", 4940000000),
    @("make", 9210000000),
    @("leave", 2600000000),
    @("tesla", 356000000),
    @("textbook", 243000000),
    @("perfect", 3990000000),
    @("usually", 1670000000),
    @("older", 1820000000),
    @("wish", 2390000000),
    @("halls", 126000000),
    @("room", 5120000000),
    @("open", 7540000000),
    @("question", 4310000000),
    @("technology", 5800000000),
    @("tesla", 365000000),
    @("serve", 3070000000),
    @("open", 9190000000),
    @("boat", 1530000000),
    @("time", 15030000000),
    @("parent", 729000000),
    @("cultural", 1180000000),
    @("textbook", 241000000),
    @("usually", 1690000000),
    @("young", 4340000000),
    @("aspect", 496000000),
    @("tesla", 370000000),
    @("perfect", 4080000000),
    @("aspect", 461000000),
    @("solve", 1880000000),
    @("answer", 2310000000),
    @("textbook", 241000000),
    @("mridangam", 868000),
    @("play", 7810000000),
    @("echo", 3310000000),
    @("older", 1150000000),
    @("shoot", 494000000),
    @("tesla", 361000000),
    @("young", 4340000000),
    @("hand", 4860000000),
    @("room", 4520000000),
    @("older", 1470000000),
    @("solve", 2510000000),
    @("wish", 2580000000),
    @("aspect", 461000000),
    @("longer", 1380000000),
    @("open", 8260000000),
    @("years", 10160000000),
    @("tesla", 361000000),
    @("longer", 1380000000),
    @("simple", 4180000000),
    @("usually", 1730000000),
    @("make", 9050000000),
    @("complete", 3400000000),
    @("parent", 774000000),
    @("response", 2170000000),
    @("science", 2960000000),
    @("limit", 978000000),
    @("curriculum", 2900000000),
    @("tesla", 375000000)
)

$startRow = 92
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
